$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Ranking swaps: two countries traded places in the table; only the
# --- displayed name moves between the two rows, row position (and thus
# --- any surrounding formatting) stays put.
$ws.Cells.Item(11,1).Value = "Sudafrica"
$ws.Cells.Item(12,1).Value = "España"

$ws.Cells.Item(55,1).Value = "Costa Rica"
$ws.Cells.Item(56,1).Value = "Nepal"

$ws.Cells.Item(204,1).Value = "Timor Oriental"
$ws.Cells.Item(205,1).Value = "Santa Lucia"

$ws.Cells.Item(214,1).Value = "Islas Malvinas"
$ws.Cells.Item(215,1).Value = "Montserrat"

# --- Refreshed statistics (Casos totales, Nuevos casos, Casos activos,
# --- Recuperados, Casos criticos, Muertes hoy, Muertes) per row.
$ws.Cells.Item(4,2).Value = 6958427
$ws.Cells.Item(4,3).Value = 32486
$ws.Cells.Item(4,4).Value = 4208887
$ws.Cells.Item(4,5).Value = 2545861
$ws.Cells.Item(4,7).Value = 512
$ws.Cells.Item(4,8).Value = 203679

$ws.Cells.Item(5,2).Value = 5392645
$ws.Cells.Item(5,3).Value = 87170
$ws.Cells.Item(5,4).Value = 4295946
$ws.Cells.Item(5,5).Value = 1009995
$ws.Cells.Item(5,7).Value = 1079
$ws.Cells.Item(5,8).Value = 86704

$ws.Cells.Item(11,2).Value = 659656
$ws.Cells.Item(11,3).Value = 2029
$ws.Cells.Item(11,4).Value = 589434
$ws.Cells.Item(11,5).Value = 54282
$ws.Cells.Item(11,7).Value = 83
$ws.Cells.Item(11,8).Value = 15940

$ws.Cells.Item(12,2).Value = 659334
$ws.Cells.Item(12,4).Value = 0
$ws.Cells.Item(12,5).Value = 0
$ws.Cells.Item(12,8).Value = 30495

$ws.Cells.Item(25,2).Value = 272080
$ws.Cells.Item(25,3).Value = 836
$ws.Cells.Item(25,5).Value = 19614

$ws.Cells.Item(29,2).Value = 142763
$ws.Cells.Item(29,3).Value = 852
$ws.Cells.Item(29,4).Value = 124184
$ws.Cells.Item(29,5).Value = 9368

$ws.Cells.Item(55,2).Value = 63712
$ws.Cells.Item(55,3).Value = 1338
$ws.Cells.Item(55,4).Value = 23552
$ws.Cells.Item(55,5).Value = 39454
$ws.Cells.Item(55,7).Value = 20
$ws.Cells.Item(55,8).Value = 706

$ws.Cells.Item(56,2).Value = 62797
$ws.Cells.Item(56,3).Value = 1204
$ws.Cells.Item(56,4).Value = 45267
$ws.Cells.Item(56,5).Value = 17129
$ws.Cells.Item(56,7).Value = 11
$ws.Cells.Item(56,8).Value = 401

$ws.Cells.Item(59,2).Value = 50992
$ws.Cells.Item(59,3).Value = 739
$ws.Cells.Item(59,4).Value = 47271
$ws.Cells.Item(59,5).Value = 3294
$ws.Cells.Item(59,7).Value = 8
$ws.Cells.Item(59,8).Value = 427

$ws.Cells.Item(99,2).Value = 9692
$ws.Cells.Item(99,3).Value = 33
$ws.Cells.Item(99,4).Value = 9341
$ws.Cells.Item(99,5).Value = 286

$ws.Cells.Item(109,2).Value = 7365
$ws.Cells.Item(109,3).Value = 4
$ws.Cells.Item(109,4).Value = 6927
$ws.Cells.Item(109,5).Value = 277

$ws.Cells.Item(126,2).Value = 4689
$ws.Cells.Item(126,3).Value = 18
$ws.Cells.Item(126,4).Value = 2910
$ws.Cells.Item(126,5).Value = 1753
$ws.Cells.Item(126,7).Value = 1
$ws.Cells.Item(126,8).Value = 26

$ws.Cells.Item(129,2).Value = 3901
$ws.Cells.Item(129,3).Value = 53
$ws.Cells.Item(129,4).Value = 1445
$ws.Cells.Item(129,5).Value = 2309

$ws.Cells.Item(143,2).Value = 3006
$ws.Cells.Item(143,3).Value = 15
$ws.Cells.Item(143,4).Value = 2349
$ws.Cells.Item(143,5).Value = 529

$ws.Cells.Item(159,2).Value = 1590
$ws.Cells.Item(159,3).Value = 25
$ws.Cells.Item(159,5).Value = 286

$ws.Cells.Item(166,2).Value = 1149
$ws.Cells.Item(166,3).Value = 2
$ws.Cells.Item(166,5).Value = 102

$ws.Cells.Item(214,4).Value = 13
$ws.Cells.Item(214,8).Value = 0

$ws.Cells.Item(215,4).Value = 12
$ws.Cells.Item(215,8).Value = 1

# --- Timestamp footer row.
$ws.Range("A1").Value = "Datos actualizados a 19 de Septiembre de 2020 a las 22:09"
